$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "393"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "891828.79"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "822"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2448771.47"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "525"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1476121.69"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "141"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334569.00"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "290"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "872172.78"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "144"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370400.26"

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "34"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "86956.00"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "35"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152929.92"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "44"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167599.15"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "194"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "493622.74"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "419"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1266532.68"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "292"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "843113.55"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "18"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70220.65"

$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "1805"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3746000.00"

$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "2787"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7335449.03"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "2845"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "6857516.41"

$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "310"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "716618.84"

$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "774"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "2324332.47"

$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value = "445"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = "1258792.91"

